# Append the latest EUR -> ARS quotation as a new row at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a leading apostrophe so Excel stores these as literal text (matching the
# existing rows) instead of auto-converting the date-/time-looking strings
# into date/time serial numbers. Clearing the formats afterwards drops the
# "quote prefix" style flag so the new cells end up with the same (default)
# style as the rest of the sheet.
$ws.Range("A90").Value = "'2025-10-20"
$ws.Range("B90").Value = "'21:22:11"
$ws.Range("C90").Value = "1.00 EUR = 1,756.2972"
$ws.Range("A90:B90").ClearFormats()
